# Slide 2 ("SEAM" bullet list): split the 4th bullet paragraph's single run
# into two runs with updated wording.
#   old: "直到每个队列都具有消息，则分别选择队列里时间戳最大的作为基础消息"
#   new: "直到每个队列都具有消息，则选择" + "所有队列里时间戳最大的作为基础消息"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(4)

# Rewrite the tail first (keeps the head's character offsets stable),
# then the head, so each Characters() call lands on the original text.
$tail = $para.Characters(16, 17)
$tail.Text = "所有队列里时间戳最大的作为基础消息"

$head = $para.Characters(1, 15)
$head.Text = "直到每个队列都具有消息，则选择"
